$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) with the new column order/labels
$ws.Cells.Item(1, 1).Value = "kitchens_1"
$ws.Cells.Item(1, 2).Value = "kitchens_2"
$ws.Cells.Item(1, 3).Value = "living_rooms_1"
$ws.Cells.Item(1, 4).Value = "bedrooms_1"
$ws.Cells.Item(1, 5).Value = "bedrooms_2"
$ws.Cells.Item(1, 6).Value = "living_rooms_2"

# Update data rows 2-7 with the new values
$data = @(
    @(0,0,0,0,1,0),
    @(0,1,0,0,0,0),
    @(0,0,0,0,0,1),
    @(0,0,0,1,0,0),
    @(1,0,0,0,0,0),
    @(0,0,1,0,0,0)
)

for ($r = 0; $r -lt 6; $r++) {
    for ($c = 0; $c -lt 6; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $data[$r][$c]
    }
}
